$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = '[Alireza%Abrishami%NULL%1,                   Shiva%Samavat%NULL%1,                   Behdad%Behnam%NULL%1,                   Mehran%Arab-Ahmadi%NULL%1,                   Mohsen%Nafar%NULL%1,                   Morteza%Sanei Taheri%NULL%1]'
$ws.Range("E3").Value2 = '[Federico%Alberici%NULL%2,                   Elisa%Delbarba%NULL%2,                   Chiara%Manenti%NULL%2,                   Laura%Econimo%NULL%2,                   Francesca%Valerio%NULL%2,                   Alessandra%Pola%NULL%2,                   Camilla%Maffei%NULL%2,                   Stefano%Possenti%NULL%2,                   Nicole%Zambetti%NULL%1,                   Marianna%Moscato%NULL%1,                   Margherita%Venturini%NULL%1,                   Stefania%Affatato%NULL%2,                   Mario%Gaggiotti%NULL%2,                   Nicola%Bossini%NULL%2,                   Francesco%Scolari%NULL%2]'
$ws.Range("E4").Value2 = '[Federico%Alberici%federico.alberici@unibs.it%0,                   Elisa%Delbarba%NULL%0,                   Chiara%Manenti%NULL%0,                   Laura%Econimo%NULL%0,                   Francesca%Valerio%NULL%0,                   Alessandra%Pola%NULL%0,                   Camilla%Maffei%NULL%0,                   Stefano%Possenti%NULL%0,                   Bernardo%Lucca%NULL%1,                   Roberta%Cortinovis%NULL%1,                   Vincenzo%Terlizzi%NULL%1,                   Mattia%Zappa%NULL%1,                   Chiara%Saccà%NULL%1,                   Elena%Pezzini%NULL%1,                   Eleonora%Calcaterra%NULL%1,                   Paola%Piarulli%NULL%1,                   Alice%Guerini%NULL%1,                   Francesca%Boni%NULL%1,                   Agnese%Gallico%NULL%1,                   Alberto%Mucchetti%NULL%1,                   Stefania%Affatato%NULL%0,                   Sergio%Bove%NULL%1,                   Martina%Bracchi%NULL%1,                   Ester Maria%Costantino%NULL%1,                   Roberto%Zubani%NULL%1,                   Corrado%Camerini%NULL%1,                   Paola%Gaggia%NULL%1,                   Ezio%Movilli%NULL%1,                   Nicola%Bossini%NULL%0,                   Mario%Gaggiotti%NULL%0,                   Francesco%Scolari%NULL%0]'
$ws.Range("E5").Value2 = '[Helena%Barrasa%NULL%0,                   Jordi%Rello%NULL%1,                   Sofia%Tejada%NULL%1,                   Alejandro%Martín%NULL%1,                   Goiatz%Balziskueta%NULL%1,                   Cristina%Vinuesa%NULL%1,                   Borja%Fernández-Miret%NULL%1,                   Ana%Villagra%NULL%1,                   Ana%Vallejo%NULL%1,                   Ana%San Sebastián%NULL%1,                   Sara%Cabañes%NULL%1,                   Sebastián%Iribarren%NULL%1,                   Fernando%Fonseca%NULL%1,                   Javier%Maynar%NULL%1,                   NULL%NULL%NULL%0]'
$ws.Range("E6").Value2 = '[Frederick S%Buckner%fbuckner@uw.edu%0,                   Denise J%McCulloch%NULL%2,                   Denise J%McCulloch%NULL%0,                   Vidya%Atluri%NULL%1,                   Michela%Blain%NULL%1,                   Sarah A%McGuffin%NULL%1,                   Arun K%Nalla%NULL%1,                   Meei-Li%Huang%NULL%1,                   Alex L%Greninger%NULL%1,                   Keith R%Jerome%NULL%1,                   Seth A%Cohen%NULL%1,                   Santiago%Neme%NULL%1,                   Margaret L%Green%NULL%1,                   Helen Y%Chu%NULL%1,                   H Nina%Kim%NULL%1]'
$ws.Range("E7").Value2 = '[Niccolò%Buetti%NULL%1,                   Timothy%Mazzuchelli%NULL%1,                   Elia%Lo Priore%NULL%1,                   Carlo%Balmelli%NULL%1,                   Michael%Llamas%NULL%1,                   Micol%Pallanza%NULL%1,                   Luigia%Elzi%NULL%1,                   Vera%Consonni%NULL%1,                   Pierpaolo%Trimboli%NULL%1,                   Valentina%Forni-Ogna%NULL%1,                   Enos%Bernasconi%NULL%1]'
$ws.Range("E8").Value2 = '[Chai%X%coreGivesNoEmail%0,                 Deliang%Huang%coreGivesNoEmail%0,                 Hong%Yu%coreGivesNoEmail%0,                 Jun%Chen%coreGivesNoEmail%0,                 Lei%Liu%coreGivesNoEmail%0,                 National%Health Commission of the People\u2019s Republic of China%coreGivesNoEmail%0,                 Pengcheng%Ou%coreGivesNoEmail%0,                 Qing%He%coreGivesNoEmail%0,                 Qingxian%Cai%coreGivesNoEmail%0,                 Rossaint%J%coreGivesNoEmail%0,                 Yang%Fu%coreGivesNoEmail%0,                 Yiming%Zhang%coreGivesNoEmail%0,                 Yinan%Su%coreGivesNoEmail%0,                 Zhang%Xia%coreGivesNoEmail%0,                 Zhenghua%Ma%coreGivesNoEmail%0,                 Zhibin%Zhu%coreGivesNoEmail%0,                 Zhiwei%Li%coreGivesNoEmail%0]'
$ws.Range("E9").Value2 = '[Bin%Cao%NULL%4,                   Yeming%Wang%NULL%0,                   Danning%Wen%NULL%1,                   Wen%Liu%NULL%1,                   Jingli%Wang%NULL%0,                   Guohui%Fan%NULL%0,                   Lianguo%Ruan%NULL%1,                   Bin%Song%NULL%0,                   Yanping%Cai%NULL%1,                   Ming%Wei%NULL%1,                   Xingwang%Li%NULL%0,                   Jiaan%Xia%NULL%0,                   Nanshan%Chen%NULL%1,                   Jie%Xiang%NULL%0,                   Ting%Yu%NULL%0,                   Tao%Bai%NULL%0,                   Xuelei%Xie%NULL%0,                   Li%Zhang%NULL%0,                   Caihong%Li%NULL%1,                   Ye%Yuan%NULL%1,                   Hua%Chen%NULL%0,                   Huadong%Li%NULL%1,                   Hanping%Huang%NULL%1,                   Shengjing%Tu%NULL%1,                   Fengyun%Gong%NULL%0,                   Ying%Liu%NULL%0,                   Yuan%Wei%NULL%0,                   Chongya%Dong%NULL%1,                   Fei%Zhou%NULL%1,                   Xiaoying%Gu%NULL%0,                   Jiuyang%Xu%NULL%0,                   Zhibo%Liu%NULL%0,                   Yi%Zhang%NULL%0,                   Hui%Li%NULL%0,                   Lianhan%Shang%NULL%1,                   Ke%Wang%NULL%0,                   Kunxia%Li%NULL%1,                   Xia%Zhou%NULL%1,                   Xuan%Dong%NULL%0,                   Zhaohui%Qu%NULL%1,                   Sixia%Lu%NULL%1,                   Xujuan%Hu%NULL%1,                   Shunan%Ruan%NULL%1,                   Shanshan%Luo%NULL%1,                   Jing%Wu%NULL%0,                   Lu%Peng%NULL%1,                   Fang%Cheng%NULL%0,                   Lihong%Pan%NULL%1,                   Jun%Zou%NULL%1,                   Chunmin%Jia%NULL%1,                   Juan%Wang%NULL%1,                   Xia%Liu%NULL%1,                   Shuzhen%Wang%NULL%1,                   Xudong%Wu%NULL%0,                   Qin%Ge%NULL%1,                   Jing%He%NULL%1,                   Haiyan%Zhan%NULL%1,                   Fang%Qiu%NULL%1,                   Li%Guo%NULL%0,                   Chaolin%Huang%NULL%1,                   Thomas%Jaki%NULL%1,                   Frederick G.%Hayden%NULL%1,                   Peter W.%Horby%NULL%1,                   Dingyu%Zhang%NULL%1,                   Chen%Wang%NULL%0]'
$ws.Range("E10").Value2 = '[Jianlei%Cao%NULL%0,                   Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%0,                   Wenlin%Cheng%NULL%0,                   Lei%Yu%NULL%0,                   Ya-Kun%Liu%NULL%0,                   Xiaoyong%Hu%NULL%0,                   Qiang%Liu%NULL%0]'
$ws.Range("E11").Value2 = '[Yang%Cao%NULL%1,                   Jia%Wei%NULL%1,                   Liang%Zou%NULL%1,                   Tiebin%Jiang%NULL%1,                   Gaoxiang%Wang%NULL%1,                   Liting%Chen%NULL%1,                   Liang%Huang%NULL%1,                   Fankai%Meng%NULL%1,                   Lifang%Huang%NULL%1,                   Na%Wang%NULL%1,                   Xiaoxi%Zhou%NULL%1,                   Hui%Luo%NULL%1,                   Zekai%Mao%NULL%1,                   Xing%Chen%NULL%1,                   Jungang%Xie%NULL%0,                   Jing%Liu%NULL%0,                   Hui%Cheng%NULL%1,                   Jianping%Zhao%NULL%0,                   Gang%Huang%NULL%1,                   Wei%Wang%NULL%0,                   Jianfeng%Zhou%NULL%1]'
$ws.Range("E12").Value2 = '[Giulio%Cavalli%cavalli.giulio@hsr.it%1,                   Giacomo%De Luca%NULL%1,                   Corrado%Campochiaro%NULL%1,                   Emanuel%Della-Torre%NULL%1,                   Marco%Ripa%NULL%1,                   Diana%Canetti%NULL%1,                   Chiara%Oltolini%NULL%1,                   Barbara%Castiglioni%NULL%1,                   Chiara%Tassan Din%NULL%1,                   Nicola%Boffini%NULL%1,                   Alessandro%Tomelleri%NULL%1,                   Nicola%Farina%NULL%1,                   Annalisa%Ruggeri%NULL%0,                   Patrizia%Rovere-Querini%NULL%0,                   Giuseppe%Di Lucca%NULL%1,                   Sabina%Martinenghi%NULL%1,                   Raffaella%Scotti%NULL%1,                   Moreno%Tresoldi%NULL%0,                   Fabio%Ciceri%NULL%1,                   Giovanni%Landoni%NULL%0,                   Alberto%Zangrillo%NULL%0,                   Paolo%Scarpellini%NULL%0,                   Lorenzo%Dagna%NULL%0]'
$ws.Range("E13").Value2 = '[Jerry Y.%Chao%NULL%1,                   Kim R.%Derespina%NULL%2,                   Betsy C.%Herold%NULL%1,                   David L.%Goldman%NULL%1,                   Margaret%Aldrich%NULL%1,                   Jacqueline%Weingarten%NULL%1,                   Henry M.%Ushay%NULL%2,                   Michael D.%Cabana%NULL%1,                   Shivanand S.%Medar%NULL%2]'
$ws.Range("E14").Value2 = '[Tao%Chen%NULL%0,                   Di%Wu%NULL%0,                   Huilong%Chen%NULL%0,                   Weiming%Yan%NULL%0,                   Danlei%Yang%NULL%0,                   Guang%Chen%NULL%0,                   Ke%Ma%NULL%0,                   Dong%Xu%NULL%0,                   Haijing%Yu%NULL%0,                   Hongwu%Wang%NULL%0,                   Tao%Wang%NULL%0,                   Wei%Guo%NULL%0,                   Jia%Chen%NULL%0,                   Chen%Ding%NULL%0,                   Xiaoping%Zhang%NULL%0,                   Jiaquan%Huang%NULL%0,                   Meifang%Han%NULL%0,                   Shusheng%Li%NULL%0,                   Xiaoping%Luo%NULL%0,                   Jianping%Zhao%NULL%0,                   Qin%Ning%NULL%0]'
$ws.Range("E15").Value2 = '[Guang%Chen%xref no email%0,           Di%Wu%xref no email%0,           Wei%Guo%xref no email%0,           Yong%Cao%xref no email%0,           Da%Huang%xref no email%0,           Hongwu%Wang%xref no email%0,           Tao%Wang%xref no email%0,           Xiaoyun%Zhang%xref no email%0,           Huilong%Chen%xref no email%0,           Haijing%Yu%xref no email%0,           Xiaoping%Zhang%xref no email%0,           Minxia%Zhang%xref no email%0,           Shiji%Wu%xref no email%0,           Jianxin%Song%xref no email%0,           Tao%Chen%xref no email%0,           Meifang%Han%xref no email%0,           Shusheng%Li%xref no email%0,           Xiaoping%Luo%xref no email%0,           Jianping%Zhao%xref no email%0,           Qin%Ning%xref no email%0]'
$ws.Range("E16").Value2 = '[Juan%Chen%NULL%1,                   Zhen-Zhen%Zhang%NULL%1,                   Yao-Kai%Chen%NULL%1,                   Quan-Xin%Long%NULL%1,                   Wen-Guang%Tian%NULL%1,                   Hai-Jun%Deng%NULL%1,                   Jie-Li%Hu%NULL%1,                   Xian-Xiang%Zhang%NULL%1,                   NULL%Pu-Liao%NULL%1,                   Jiang-Lin%Xiang%NULL%1,                   Dao-Xin%Wang%NULL%1,                   Peng%Hu%NULL%0,                   Fa-Chun%Zhou%NULL%1,                   Zhi-Jie%Li%NULL%1,                   Hong-Mei%Xu%NULL%1,                   Xue-Fei%Cai%NULL%1,                   De-Qiang%Wang%NULL%1,                   Yuan%Hu%NULL%1,                   Ni%Tang%NULL%1,                   Bei-Zhong%Liu%liubeizhong@cqmu.edu.cn%1,                   Gui-Cheng%Wu%wuguic@hotmail.com%1,                   Ai-Long%Huang%ahuang@cqmu.edu.cn%1]'
$ws.Range("E17").Value2 = '[Nanshan%Chen%NULL%0,                   Min%Zhou%NULL%0,                   Xuan%Dong%NULL%0,                   Jieming%Qu%NULL%0,                   Fengyun%Gong%NULL%0,                   Yang%Han%NULL%0,                   Yang%Qiu%NULL%0,                   Jingli%Wang%NULL%0,                   Ying%Liu%NULL%0,                   Yuan%Wei%NULL%0,                   Jia''an%Xia%NULL%0,                   Ting%Yu%NULL%0,                   Xinxin%Zhang%NULL%0,                   Li%Zhang%NULL%0]'
$ws.Range("E18").Value2 = '[Qingqing%Chen%NULL%0,                   Zhencang%Zheng%NULL%1,                   Chao%Zhang%NULL%1,                   Xijiang%Zhang%NULL%1,                   Huijuan%Wu%NULL%1,                   Jingdong%Wang%NULL%1,                   Shuwei%Wang%NULL%1,                   Cheng%Zheng%dr.zhengcheng@foxmail.com%1]'
$ws.Range("E19").Value2 = '[Ruchong%Chen%NULL%0,                   Wenhua%Liang%NULL%0,                   Mei%Jiang%NULL%0,                   Weijie%Guan%NULL%0,                   Chen%Zhan%NULL%0,                   Tao%Wang%NULL%0,                   Chunli%Tang%NULL%0,                   Ling%Sang%NULL%0,                   Jiaxing%Liu%NULL%0,                   Zhengyi%Ni%NULL%0,                   Yu%Hu%NULL%0,                   Lei%Liu%NULL%0,                   Hong%Shan%NULL%0,                   Chunliang%Lei%NULL%0,                   Yixiang%Peng%NULL%0,                   Li%Wei%NULL%0,                   Yong%Liu%NULL%0,                   Yahua%Hu%NULL%0,                   Peng%Peng%NULL%0,                   Jianming%Wang%NULL%0,                   Jiyang%Liu%NULL%0,                   Zhong%Chen%NULL%0,                   Gang%Li%NULL%0,                   Zhijian%Zheng%NULL%0,                   Shaoqin%Qiu%NULL%0,                   Jie%Luo%NULL%0,                   Changjiang%Ye%NULL%0,                   Shaoyong%Zhu%NULL%0,                   Xiaoqing%Liu%NULL%0,                   Linling%Cheng%NULL%0,                   Feng%Ye%NULL%0,                   Jinping%Zheng%NULL%0,                   Nuofu%Zhang%NULL%0,                   Yimin%Li%NULL%0,                   Jianxing%He%NULL%0,                   Shiyue%Li%lishiyue@188.com%0,                   Nanshan%Zhong%NULL%0,                   NULL%NULL%NULL%0]'
$ws.Range("E20").Value2 = '[Min%Chen%NULL%1,                   Wei%An%NULL%2,                   Wei%An%NULL%0,                   Fei%Xia%NULL%0,                   Ping%Yang%NULL%1,                   Kuangyu%Li%NULL%1,                   Qin%Zhou%NULL%1,                   Shasha%Fang%NULL%1,                   Yaling%Liao%NULL%1,                   Xin%Xu%NULL%0,                   Jialin%Liu%NULL%1,                   Shiguo%Liu%NULL%1,                   Tao%Qin%NULL%1,                   Jianjun%Zhang%NULL%1,                   Wei%Wei%NULL%0,                   Yafang%Zhang%NULL%1,                   Guowei%Zhang%NULL%1,                   Mingwei%Zhang%383193522@qq.com%1]'
$ws.Range("E21").Value2 = '[Jiaojiao%Chu%NULL%0,                   Nan%Yang%NULL%2,                   Nan%Yang%NULL%0,                   Yanqiu%Wei%NULL%1,                   Huihui%Yue%NULL%1,                   Fengqin%Zhang%NULL%1,                   Jianping%Zhao%NULL%0,                   Li%He%NULL%1,                   Gaohong%Sheng%NULL%2,                   Gaohong%Sheng%NULL%0,                   Peng%Chen%NULL%1,                   Gang%Li%NULL%0,                   Sisi%Wu%NULL%1,                   Bo%Zhang%NULL%0,                   Shu%Zhang%szhang@tjh.tjmu.edu.cn%1,                   Congyi%Wang%wangcy@tjh.tjmu.edu.cn%1,                   Xiaoping%Miao%miaoxp@hust.edu.cn%1,                   Juan%Li%lijuan@tjh.tjmu.edu.cn%0,                   Wenhua%Liu%liuwh_2013@126.com%1,                   Huilan%Zhang%huilanz_76@163.com%1]'
$ws.Range("E22").Value2 = '[Rebecca%Cosgriff%NULL%1,                   Susannah%Ahern%NULL%1,                   Scott C.%Bell%NULL%1,                   Keith%Brownlee%NULL%1,                   Pierre-Régis%Burgel%NULL%1,                   Cass%Byrnes%NULL%1,                   Harriet%Corvol%NULL%1,                   Stephanie Y.%Cheng%NULL%1,                   Alexander%Elbert%NULL%1,                   Albert%Faro%NULL%1,                   Christopher H.%Goss%NULL%1,                   Vincent%Gulmans%NULL%1,                   Bruce C.%Marshall%NULL%1,                   Edward%McKone%NULL%1,                   Peter G.%Middleton%NULL%1,                   Rasa%Ruseckaite%NULL%1,                   Anne L.%Stephenson%NULL%1,                   Siobhán B%Carr%NULL%1,                   David%Reid%NULL%1,                   Peter%Wark%NULL%1,                   Géraldine%Daneau%NULL%1,                   Véronique%Boussaud%NULL%1,                   Graziella%Brinchault%NULL%1,                   Emmanuelle%Coirier-Duet%NULL%1,                   Jean-Christophe%Dubus%NULL%1,                   Dominique%Grenet%NULL%1,                   Sandra de%Miranda%NULL%1,                   Laurence%Beaumont%NULL%1,                   Reem%Kanaan%NULL%1,                   Muriel%Lauraens%NULL%1,                   Clémence%Martin%NULL%1,                   Marie%Mittaine%NULL%1,                   Anne%Prévotat%NULL%1,                   Martine%Reynaud-Gaubert%NULL%1,                   Isabelle%Sermet-Gaudelus%NULL%1,                   Aurelie%Tatopoulos%NULL%1,                   Lutz%Nährlich%NULL%1,                   Barry%Plant%NULL%1,                   Cedric%Gunaratnam%NULL%1,                   Abaigeal%Jackson%NULL%1,                   Karin M. de Winter-de%Groot%NULL%1,                   Bart%Luijk%NULL%1,                   Geertjan%Wesseling%NULL%1,                   Mark%Allenby%NULL%1,                   Jamie%Duckers%NULL%1,                   Elaine%Gunn%NULL%1,                   Andrew%Jones%NULL%1,                   Robert Ian%Ketchell%NULL%1,                   Susan L.%Madge%NULL%1,                   Anirban%Maitra%NULL%1,                   Ghulam%Mujtaba%NULL%1,                   Helen%Rodgers%NULL%1,                   Nadia%Shafi%NULL%1,                   Nicholas%Simmonds%NULL%1,                   Kevin%Southern%NULL%1,                   Danie%Watson%NULL%1,                   Samar%Rizvi%NULL%1,                   Julie%Seguin%NULL%1,                   Janet%Garbarz%NULL%1,                   Kristen%Rosamilia%NULL%1,                   Maria%Berdella%NULL%1,                   Jerry A.%Nick%NULL%1,                   Richard%Belkin%NULL%1,                   Diana%Gilmore%NULL%1,                   Kim%McBennett%NULL%1,                   Rita%Padoan%NULL%1,                   Marco%Salvatore%NULL%1]'
$ws.Range("E23").Value2 = '[Matthew J%Cummings%NULL%0,                   Matthew R%Baldwin%NULL%0,                   Darryl%Abrams%NULL%1,                   Samuel D%Jacobson%NULL%1,                   Benjamin J%Meyer%NULL%1,                   Elizabeth M%Balough%NULL%1,                   Justin G%Aaron%NULL%1,                   Jan%Claassen%NULL%0,                   LeRoy E%Rabbani%NULL%1,                   Jonathan%Hastie%NULL%1,                   Beth R%Hochman%NULL%1,                   John%Salazar-Schicchi%NULL%1,                   Natalie H%Yip%NULL%1,                   Daniel%Brodie%NULL%0,                   Max R%O''Donnell%NULL%1]'
$ws.Range("E24").Value2 = '[Hui%Dai%NULL%1,                   Xin%Zhang%NULL%0,                   Jianguo%Xia%NULL%1,                   Tao%Zhang%NULL%1,                   Yalei%Shang%NULL%1,                   Renjun%Huang%NULL%1,                   Rongrong%Liu%NULL%1,                   Dan%Wang%NULL%0,                   Min%Li%105549156@qq.com%1,                   Jinping%Wu%czwjp@sina.com%1,                   Qiuzhen%Xu%xuqiuzhen831@sina.com%1,                   Yonggang%Li%liyonggang224@163.com%1]'
